$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New config row: maSoThue / 0106777886
$ws.Range("A15").Value = "maSoThue"

# Force the tax-code value to be stored as text so the leading zero survives
# (typed with a leading apostrophe, like a user would in Excel -> quotePrefix style)
$ws.Range("B15").Value = "'0106777886"

# Move selection the way the author's session ended up (one row below the new row)
$ws.Range("B16").Select() | Out-Null
